$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 99
$ws.Range("I2").Value = 99
$ws.Range("K2").Value = 99
$ws.Range("M2").Value = 14
$ws.Range("H18").Value = 1036.2307
$ws.Range("I18").Value = 1036.2307
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1036.2307
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -752.2307000000001
$ws.Range("N18").ClearContents()
$ws.Range("H28").Value = 273.1875
$ws.Range("I28").Value = 270.20834
$ws.Range("J28").Value = 282.125
$ws.Range("K28").Value = 270.20834
$ws.Range("L28").Value = 282.125
$ws.Range("M28").Value = 214.79166
$ws.Range("N28").Value = -1252.125
$ws.Range("H32").Value = 1114.8572
$ws.Range("I32").Value = 833.3333
$ws.Range("J32").Value = 1326
$ws.Range("K32").Value = 833.3333
$ws.Range("L32").Value = 1326
$ws.Range("M32").Value = -507.3333
$ws.Range("N32").Value = -1978
$ws.Range("H64").Value = 3191.9429
$ws.Range("I64").Value = 2921.2666
$ws.Range("K64").Value = 2921.2666
$ws.Range("M64").Value = -2673.2666
$ws.Range("H67").Value = 3191.9429
$ws.Range("I67").Value = 2921.2666
$ws.Range("K67").Value = 2921.2666
$ws.Range("M67").Value = -2063.2666
$ws.Range("H93").Value = 36685.715
$ws.Range("J93").Value = 36685.715
$ws.Range("L93").Value = 36685.715
$ws.Range("N93").Value = -41677.715
$ws.Range("H95").Value = 30111
$ws.Range("J95").Value = 30111
$ws.Range("L95").Value = 30111
$ws.Range("N95").Value = -35603
$ws.Range("H116").Value = 12502778
$ws.Range("J116").Value = 3766.3333
$ws.Range("L116").Value = 3766.3333
$ws.Range("N116").Value = -10650.3333
$ws.Range("H135").Value = 1289.4286
$ws.Range("I135").Value = 1117.4375
$ws.Range("K135").Value = 10056.9375
$ws.Range("M135").Value = -7521.9375
$ws.Range("H137").Value = 1362.9269
$ws.Range("I137").Value = 1121.871
$ws.Range("J137").Value = 2110.2
$ws.Range("K137").Value = 3365.613
$ws.Range("L137").Value = 6330.599999999999
$ws.Range("M137").Value = -815.6130000000003
$ws.Range("N137").Value = -11430.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 32000
$ws.Range("J64").Value = 32000
$ws.Range("L64").Value = 32000
$ws.Range("N64").Value = -32496
$ws.Range("H67").Value = 32000
$ws.Range("J67").Value = 32000
$ws.Range("L67").Value = 32000
$ws.Range("N67").Value = -33716
$ws.Range("H74").Value = 1424.1538
$ws.Range("I74").Value = 1168.7778
$ws.Range("J74").Value = 1998.75
$ws.Range("K74").Value = 1168.7778
$ws.Range("L74").Value = 1998.75
$ws.Range("M74").Value = -294.7778000000001
$ws.Range("N74").Value = -3746.75
$ws.Range("H77").Value = 1424.1538
$ws.Range("I77").Value = 1168.7778
$ws.Range("J77").Value = 1998.75
$ws.Range("K77").Value = 5843.889
$ws.Range("L77").Value = 9993.75
$ws.Range("M77").Value = -1475.889
$ws.Range("N77").Value = -18729.75
$ws.Range("H92").Value = 66989
$ws.Range("J92").Value = 66989
$ws.Range("L92").Value = 66989
$ws.Range("N92").Value = -71981

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1532.0834
$ws.Range("I134").Value = 1411.025
$ws.Range("K134").Value = 4233.075000000001
$ws.Range("M134").Value = -1698.075000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2046.7142
$ws.Range("I31").Value = 1439.3784
$ws.Range("J31").Value = 3919.3333
$ws.Range("K31").Value = 1439.3784
$ws.Range("L31").Value = 3919.3333
$ws.Range("M31").Value = -1144.3784
$ws.Range("N31").Value = -4509.3333
$ws.Range("H34").Value = 2046.7142
$ws.Range("I34").Value = 1439.3784
$ws.Range("J34").Value = 3919.3333
$ws.Range("K34").Value = 1439.3784
$ws.Range("L34").Value = 3919.3333
$ws.Range("M34").Value = -1237.3784
$ws.Range("N34").Value = -4323.3333
$ws.Range("H58").Value = 2194.3
$ws.Range("J58").Value = 3995
$ws.Range("L58").Value = 3995
$ws.Range("N58").Value = -4401
$ws.Range("H132").Value = 1907
$ws.Range("I132").Value = 1411.2222
$ws.Range("J132").Value = 2799.4
$ws.Range("K132").Value = 4233.6666
$ws.Range("L132").Value = 8398.200000000001
$ws.Range("M132").Value = -1703.6666
$ws.Range("N132").Value = -13458.2
$ws.Range("H134").Value = 2305.75
$ws.Range("I134").Value = 2291.1
$ws.Range("K134").Value = 6873.299999999999
$ws.Range("M134").Value = -4338.299999999999
$ws.Range("H136").Value = 2194.3
$ws.Range("J136").Value = 3995
$ws.Range("L136").Value = 11985
$ws.Range("N136").Value = -17085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14926816
$ws.Range("J131").Value = 15874524
$ws.Range("L131").Value = 47623572
$ws.Range("N131").Value = -47633652
$ws.Range("H132").Value = 1167
$ws.Range("I132").Value = 1167
$ws.Range("K132").Value = 10503
$ws.Range("M132").Value = -7973

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5150
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 5780
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 5780
$ws.Range("M43").Value = -1849
$ws.Range("N43").Value = -6082
$ws.Range("H46").Value = 11250
$ws.Range("J46").Value = 11906.25
$ws.Range("L46").Value = 11906.25
$ws.Range("N46").Value = -12218.25
$ws.Range("H58").Value = 15957.143
$ws.Range("I58").Value = 9900
$ws.Range("J58").Value = 20500
$ws.Range("K58").Value = 9900
$ws.Range("L58").Value = 20500
$ws.Range("M58").Value = -9623
$ws.Range("N58").Value = -21054
$ws.Range("H92").Value = 26490
$ws.Range("J92").Value = 26490
$ws.Range("L92").Value = 26490
$ws.Range("N92").Value = -30234
$ws.Range("H95").Value = 1264479.2
$ws.Range("J95").Value = 1264479.2
$ws.Range("L95").Value = 1264479.2
$ws.Range("N95").Value = -1269971.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1840.4
$ws.Range("I16").Value = 1566.6666
$ws.Range("J16").Value = 2251
$ws.Range("K16").Value = 1566.6666
$ws.Range("L16").Value = 2251
$ws.Range("M16").Value = -1396.6666
$ws.Range("N16").Value = -2591
$ws.Range("H46").Value = 1833.3334
$ws.Range("I46").Value = 1666.6666
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1666.6666
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1478.6666
$ws.Range("N46").Value = -2376
$ws.Range("H56").Value = 27771.5
$ws.Range("I56").Value = 27333.334
$ws.Range("J56").Value = 28034.4
$ws.Range("K56").Value = 27333.334
$ws.Range("L56").Value = 28034.4
$ws.Range("M56").Value = -26642.334
$ws.Range("N56").Value = -29416.4
$ws.Range("H82").Value = 1558.9333
$ws.Range("I82").Value = 1180.091
$ws.Range("J82").Value = 2600.75
$ws.Range("K82").Value = 1180.091
$ws.Range("L82").Value = 2600.75
$ws.Range("M82").Value = -819.0909999999999
$ws.Range("N82").Value = -3322.75
$ws.Range("H85").Value = 1558.9333
$ws.Range("I85").Value = 1180.091
$ws.Range("J85").Value = 2600.75
$ws.Range("K85").Value = 1180.091
$ws.Range("L85").Value = 2600.75
$ws.Range("M85").Value = 67.90900000000011
$ws.Range("N85").Value = -5096.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 91583.27
$ws.Range("I81").Value = 91583.27
$ws.Range("K81").Value = 183166.54
$ws.Range("M81").Value = -182105.54
$ws.Range("H84").Value = 91583.27
$ws.Range("I84").Value = 91583.27
$ws.Range("K84").Value = 915832.7000000001
$ws.Range("M84").Value = -910528.7000000001
